# MENT-169: Create New Question Categories
# Adds two new category rows below the existing list in column A:
#   A14 -> Laboratório
#   A15 -> UATS

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (A13) down onto the
# two new rows so they pick up the same cell style (font "Times New
# Roman", color FF333333 - cellXfs index 1) as the rest of the list,
# instead of defaulting to the workbook's base style.
$ws.Range("A13").Copy()
$ws.Range("A14:A15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new category names.
$ws.Range("A14").Value = "Laboratório"
$ws.Range("A15").Value = "UATS"

# Match the author's final selection state.
$ws.Range("A25").Select()
